$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Rename the "ID"-suffixed headers to their shorter form (hide filter in locations lists)
$ws.Range("A1").Value = "Category"
$ws.Range("B1").Value = "Errand Type"
$ws.Range("D1").Value = "Addon1 Errand Type"
$ws.Range("F1").Value = "Addon2 Errand Type"

# Clear the sample data row, but keep H2's existing style
$ws.Range("L2:M2").Style = "Normal"
$ws.Range("A2:G2").ClearContents()
$ws.Range("H2").ClearContents()
$ws.Range("J2:M2").ClearContents()

# Move the active selection back to A2
$ws.Range("A2").Select() | Out-Null
